$wb = $excel.ActiveWorkbook

# --- 1. Add a new "Sheet1" as the first tab, holding the vendor list ---
$first = $wb.Worksheets.Item(1)
$vendorSheet = $wb.Worksheets.Add($first)

$vendors = @("Vendors", "Vendor 1", "Vendor 2", "Vendor 3", "Vendor 4", "Vendor 5", "Vendor 6")
for ($i = 0; $i -lt $vendors.Length; $i++) {
    $vendorSheet.Cells.Item($i + 1, 1).Value = $vendors[$i]
}

# --- 2. Compact each existing vendor sheet's data from columns A/C/E into A/B/C ---
$ws1 = $wb.Worksheets.Item("Vendor 1")
$ws1.Range("C1:C5").Cut($ws1.Range("B1:B5"))
$ws1.Range("E1:E5").Cut($ws1.Range("C1:C5"))
$ws1.Range("D1:E5").Clear()
$null = $ws1.Range("C1:C5").Select()

$ws2 = $wb.Worksheets.Item("Vendor 2")
$ws2.Range("C1:C5").Cut($ws2.Range("B1:B5"))
$ws2.Range("E1:E5").Cut($ws2.Range("C1:C5"))
$ws2.Range("D1:E5").Clear()
$null = $ws2.Range("C1:C5").Select()

$ws3 = $wb.Worksheets.Item("Vendor 3")
$ws3.Range("C1:C5").Cut($ws3.Range("B1:B5"))
$ws3.Range("E1:E5").Cut($ws3.Range("C1:C5"))
$ws3.Range("D1:E5").Clear()
$null = $ws3.Range("F6").Select()

# --- 3. Leave the new vendor list sheet selected/active, as in the target workbook ---
$null = $vendorSheet.Range("D4").Select()
